$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2233
$ws.Range("L3").Value = 2253
$ws.Range("E4").Value = 687
$ws.Range("F4").Value = 691
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 616
$ws.Range("L5").Value = 134
$ws.Range("L6").Value = 2030
$ws.Range("E7").Value = 8430
$ws.Range("F7").Value = 7722
$ws.Range("K7").Value = 9224
$ws.Range("L7").Value = 7266

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 58
$ws.Range("L4").Value = 29
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 236
$ws.Range("L8").Value = 460
$ws.Range("L9").Value = 47
$ws.Range("L11").Value = 130
$ws.Range("L18").Value = 52
$ws.Range("L19").Value = 208
$ws.Range("L23").Value = 74
$ws.Range("L27").Value = 74
$ws.Range("L29").Value = 372
$ws.Range("L33").Value = 327
$ws.Range("L36").Value = 104
$ws.Range("L37").Value = 261
$ws.Range("E41").Value = 54
$ws.Range("L42").Value = 226
$ws.Range("L43").Value = 58
$ws.Range("L44").Value = 52
$ws.Range("L47").Value = 56
$ws.Range("L48").Value = 99
$ws.Range("L52").Value = 144
$ws.Range("L53").Value = 92
$ws.Range("L54").Value = 147
$ws.Range("L57").Value = 31
$ws.Range("L58").Value = 4
$ws.Range("L59").Value = 11
$ws.Range("F63").Value = 85
$ws.Range("K63").Value = 44
$ws.Range("L63").Value = 24
$ws.Range("L66").Value = 16
$ws.Range("L67").Value = 258
$ws.Range("L68").Value = 20
$ws.Range("L73").Value = 57
$ws.Range("L77").Value = 44
$ws.Range("L78").Value = 98
$ws.Range("L79").Value = 201
$ws.Range("L85").Value = 380
$ws.Range("L86").Value = 56
$ws.Range("L87").Value = 21
$ws.Range("L88").Value = 102
$ws.Range("L91").Value = 103
$ws.Range("L94").Value = 88
$ws.Range("L96").Value = 70
$ws.Range("L97").Value = 68
$ws.Range("L98").Value = 52
$ws.Range("L99").Value = 115
$ws.Range("E101").Value = 8430
$ws.Range("F101").Value = 7722
$ws.Range("K101").Value = 9224
$ws.Range("L101").Value = 7266

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 69
$ws.Range("L3").Value = 75
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 45
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 117
$ws.Range("L3").Value = 158
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 380

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 131
$ws.Range("L3").Value = 155
$ws.Range("L4").Value = 34
$ws.Range("L6").Value = 121
$ws.Range("L7").Value = 460

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 105
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 114
$ws.Range("L7").Value = 327

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 76
$ws.Range("L7").Value = 261

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 28
$ws.Range("L3").Value = 52
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 86
$ws.Range("L7").Value = 258

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 121
$ws.Range("L3").Value = 133
$ws.Range("L7").Value = 372

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 68
$ws.Range("L6").Value = 65
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L3").Value = 11
$ws.Range("E4").Value = 5
$ws.Range("L4").Value = 2
$ws.Range("E7").Value = 54

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L4").Value = 23
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 32
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L4").Value = 15
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 24
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 19
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L4").Value = 4
$ws.Range("L6").Value = 16

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 33
$ws.Range("L6").Value = 56

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 9
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L2").Value = 7
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("L5").ClearContents()
$ws.Range("L6").Value = 4
